$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column D to match new, longer Results text ---
$ws.Columns("D").ColumnWidth = 41.75

function Set-Plain($addr, $text) {
    $ws.Range($addr).Value = $text
}

function Set-SpecialBC($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Value = $text
    $rng.Borders.LineStyle = 0
    $rng.Interior.Pattern = -4142
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4160
    $rng.WrapText = $false
    $rng.Font.Name = "Arial"
}

function Set-SpecialD($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.Value = $text
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4160
    $rng.WrapText = $true
}

function Set-Rich($addr, $introText, $boldText) {
    $rng = $ws.Range($addr)
    $full = $introText + $boldText
    $rng.Value = $full
    $introLen = $introText.Length
    $rng.Characters(1, $introLen).Font.Bold = $false
    $rng.Characters($introLen + 1, $boldText.Length).Font.Bold = $true
}

# --- Row 3 (Test 1): add new Results text ---
Set-Plain "D3" "Menu displays correctly, if database_data does not exist, creates the file"
Set-Plain "D4" "displays `"database is empty`""
Set-SpecialBC "B5" "Enter an invalid input"
Set-SpecialBC "C5" "The program displays a error message"
Set-SpecialD "D5" "displays `"error: option does not exists"
Set-Rich "B6" "Select the Add Option and add the following information:`n`n" "Account#: 1000`nName: test name 1`nAddress: test address 1"
Set-Plain "C6" "The new record is added and confirmation is displayed to the user. Use printall to confirm"
Set-Plain "D6" "displays `"The account was added succesfully`", printall correctly displays address"
Set-Rich "B7" "Select the Add Option and add the following information:`n`n" "Account#: 1002`nName: test name 3`nAddress: test address 3"
Set-Plain "C7" "The new record is added and confirmation is displayed to the user. Use printall to confirm"
Set-Plain "D7" "displays `"The account was added succesfully`", printall correctly displays addresses in correct order"
Set-Rich "B8" "Select the Add Option and add the following information:`n`n" "Account#: 1001`nName: test name 2`nAddress: test address 2"
Set-Plain "C8" "The new record is added and confirmation is displayed to the user. Use printall to confirm"
Set-Plain "D8" "displays `"The account was added succesfully`", printall correctly displays addresses in correct order"
Set-Plain "A9" "Test 7"
Set-Plain "B9" "Select Find Option and input Account#: 1000"
Set-Plain "C9" "The record is found and confirmation is displayed to the user. "
Set-Plain "D9" "The correct record is found and printed"
Set-Plain "A10" "Test 8"
Set-Plain "B10" "Select Delete Option and delete 1000"
Set-Plain "C10" "The record is deleted and confirmation is displayed to the user. "
Set-Plain "D10" "The delete option DOES NOT perform correctly, unable to properly delete function"
Set-Plain "A11" " Test 9"
Set-Plain "B11" "Repeat Test 7"
Set-Plain "C11" "The record is not found and a error is displayed to the user. "
Set-Plain "D11" "The delete did not perform correctly, unable to test"
Set-Plain "A12" "Test 10"
Set-Plain "B12" "Repeat Test 2"
Set-Plain "C12" "The program displays all records left (1002, 1001)."
Set-Plain "D12" "The program correctly displays the remaining 3 since delete did not function correctly"
Set-Plain "A13" "Test 11"
Set-Plain "B13" "Select the Quit Option"
Set-Plain "C13" "The program quits"
Set-Plain "D13" "The program displays `"Exiting the program...`" and succesfully quits"
